$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column E to fit the new longer text
$ws.Columns.Item(5).ColumnWidth = 54.28515625

# New shared strings must be created in this order so the shared string
# table layout matches: B39, B40, E39, E40
$entryDate = Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0

# Row 39: 2149. Rearrange Array Elements by Sign (Java) - Varitey 1
$ws.Cells.Item(39, 1).Value = 2149
$ws.Cells.Item(39, 2).Value = "Rearrange Array Elements by Sign"

# Row 40: GFG. Alternate positive and negative numbers (Java) - Varitey 2
$ws.Cells.Item(40, 1).Value = "GFG"
$ws.Cells.Item(40, 2).Value = "Alternate positive and negative numbers"

$ws.Cells.Item(39, 5).Value = "Varitey 1- where +ve and -ve are same in number"
$ws.Cells.Item(40, 5).Value = "Varitey 2- where +ve and -ve are not in number"

$ws.Cells.Item(39, 3).Value = "Java"
$ws.Cells.Item(39, 4).Value = $entryDate
$ws.Cells.Item(40, 3).Value = "Java"
$ws.Cells.Item(40, 4).Value = $entryDate

# Apply formatting to match the style used for row 38 (same block as rows 35-38)
$srcRange = $ws.Range("A38:E38")
$destRange = $ws.Range("A39:E40")
$srcRange.Copy()
$destRange.PasteSpecial(-4122) # xlPasteFormats

# Re-set the values since PasteSpecial formats only, values set above stand
$ws.Cells.Item(39, 4).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(40, 4).NumberFormat = "m/d/yyyy"

# Update selection / view
$ws.Range("B44").Select()
